# Concertina Wave Calculator - update input parameters
#
# A2 (wave amplitude / joint angle input) changes from 45 -> 30.
# D2 used to be a formula (=2.67*180/PI()) that produced ~153, which fed
# back through the geometry chain to an out-of-range starting heading;
# the refactor replaces it with the "as-used" literal value so the
# concertina motion calc starts from the correct heading.
# E2 (initial heading input) changes from 15 -> 13.5.
# All other changed cells in the sheet (and the chart caches that mirror
# them) are formulas that recompute automatically from these three inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 30
$ws.Range("D2").Value = 47.664416515956525
$ws.Range("E2").Value = 13.5

# Restore the view: scroll back to the top-left of the sheet and move the
# active selection from D3 to D2 (matches the sheetView/selection diff).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D2").Select()

$wb.Application.CalculateFull()

# Make sure every chart embedded on the sheet picks up the recalculated
# series data (the charts plot J/K/L columns that all depend on A2/D2/E2).
$charts = $ws.ChartObjects()
for ($i = 1; $i -le $charts.Count(); $i++) {
    $co = $charts.Item($i)
    $co.Chart.Refresh()
}
